# Insert a new data row at row 45 (pushing existing rows 45..168 down to 46..169)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 45; existing rows 45+ shift down to 46+.
$ws.Rows.Item(45).Insert()

# Fill in the new row 45 with the new weekly record's values.
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(45, 3).Value = 'La Araucanía'
$ws.Cells.Item(45, 4).Value = 44525
$ws.Cells.Item(45, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(45, 5).Value = 9
$ws.Cells.Item(45, 6).Value = 'Fruta'
$ws.Cells.Item(45, 7).Value = 100102
$ws.Cells.Item(45, 8).Value = 'Cítricos'
$ws.Cells.Item(45, 9).Value = 100102006
$ws.Cells.Item(45, 10).Value = 'Pomelo'
$ws.Cells.Item(45, 11).Value = 'Start Ruby'
$ws.Cells.Item(45, 12).Value = 'Primera'
$ws.Cells.Item(45, 13).Value = 25
$ws.Cells.Item(45, 14).Value = 15000
$ws.Cells.Item(45, 15).Value = 15000
$ws.Cells.Item(45, 16).Value = 15000
$ws.Cells.Item(45, 17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(45, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(45, 19).Value = 1071
$ws.Cells.Item(45, 20).Value = 14
